# The deck's single slide master/theme ("Integral") is swapped back to the
# default "Office Theme" color palette (the companion theme part in the
# package, originally only wired to the Notes Master, held the Office
# Theme colors - this macro re-applies that palette onto the live design so
# the deck's effective theme colors match the Office Theme again).

function ToRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$tcs.Item(1).RGB  = ToRgb 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = ToRgb 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = ToRgb 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = ToRgb 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = ToRgb 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = ToRgb 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = ToRgb 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = ToRgb 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = ToRgb 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = ToRgb 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = ToRgb 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = ToRgb 0x95 0x4F 0x72   # folHlink
